$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 104: date 2021-07-22 (44399) -> 2021-09-09 (44448), volume 100 -> 55
$ws.Range("D104").Value = 44448
$ws.Range("J104").Value = 55

# Update existing row 105: date 2021-07-23 (44400) -> 2021-07-22 (44399), volume 40 -> 100
$ws.Range("D105").Value = 44399
$ws.Range("J105").Value = 100

# Insert new row 106 - a copy of what row 105 used to be before this edit
$ws.Range("A106").Value = 10
$ws.Range("B106").Value = "Vega Modelo de Temuco"
$ws.Range("C106").Value = "La Araucanía"
$ws.Range("D106").Value = 44400
$ws.Range("D106").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E106").Value = 9
$ws.Range("F106").Value = 100114007
$ws.Range("G106").Value = "Jengibre"
$ws.Range("H106").Value = "Sin especificar"
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 40
$ws.Range("K106").Value = 20000
$ws.Range("L106").Value = 20000
$ws.Range("M106").Value = 20000
$ws.Range("N106").Value = "$/caja 13 kilos"
$ws.Range("O106").Value = "Perú"
$ws.Range("P106").Value = 1538
$ws.Range("Q106").Value = 13
$ws.Range("R106").Value = "Hortaliza"
